# Week_06.docx - "added message about pw"
#
# Insert a new bold bullet-list item warning students not to create a
# second account if they forget their password, placed immediately
# before the existing "TICKABLE It is very easy to get help in Sage..."
# bullet (i.e. right after the "Video hint" bullet that precedes it).

$d = $word.ActiveDocument

# Locate the unique paragraph that starts the "It is very easy to get
# help in Sage" TICKABLE item - this anchors us precisely even though
# "TICKABLE" / "Video hint" text recurs throughout the document.
$anchor = $d.Content
$found = $anchor.Find.Execute("It is very easy to get help in Sage", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'It is very easy to get help in Sage' paragraph"
}

$targetPara = $anchor.Paragraphs(1)
$precedingPara = $targetPara.Previous()

# Insert a brand-new empty paragraph right after the preceding ("Video
# hint") paragraph - it inherits that paragraph's list (numId 1) and
# carries no leftover character styling, so we can fill it in cleanly.
$precedingPara.Range.InsertParagraphAfter()
$newPara = $precedingPara.Next()

$newPara.Range.Text = "If you forget your password DO NOT CREATE ANOTHER ACCOUNT: come and see me (Vince Knight) and I can reset your password."
$newPara.Range.Font.Bold = $true

Write-Host "Inserted password-reset paragraph:" $newPara.Range.Text
Write-Host "Followed by:" $newPara.Next().Range.Text
